# TC32_Canine_Filter_Breed-IrishWlfhnd.xlsx  -  "Fixed ICDC breed all testcases"
#
# The StatQuery column (C) on the "startup" sheet held a stale/broken Cypher
# query used to compute summary statistics. Replace it with the corrected
# query (now returning Programs/Studies/Cases/Samples/Case Files/Study Files)
# on every data row (CasesTab, SamplesTab, FilesTab). Also refresh the sheet
# view (zoom + selection) to match the author's final save state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Irish Wolfhound']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Replace the StatQuery text (column C) for each of the 3 data rows
# (CasesTab / SamplesTab / FilesTab) with the corrected query.
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Refresh the view: zoomed to 100% and selection sitting on B4, matching the
# workbook as last saved.
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("B4").Select()
